# Tuntikirjanpito.xlsx - add three new work-log entries (Hook/Redux testing
# study, jest setup, PrivateRoute tests) below the last existing entry and
# move the totals block down to make room, updating the SUM/percentage
# formulas and the sheet's used-range/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- make room ----------------------------------------------------------
# Row 114 used to be a single blank spacer row directly above the totals
# block (rows 115-117). Insert 3 fresh rows there for the new log entries,
# then one more blank spacer row so the totals block ends up two rows
# below the last entry again, just wider than before.
$ws.Rows("114:116").Insert() | Out-Null
$ws.Rows("118").Insert() | Out-Null

# --- new entry: row 114 (no date cell) ----------------------------------
$ws.Range("A114").Clear() | Out-Null
$ws.Range("B114").Value = 1.5
$ws.Range("C114").Value = "Hook (useSelector, useDispatch) testien selvittelyä ja opiskelua"
$ws.Range("D114").Value = "client"

# --- new entry: row 115 (1.2.2022) --------------------------------------
$ws.Range("A113").Copy() | Out-Null
$ws.Range("A115").PasteSpecial(-4122) | Out-Null
$ws.Range("A115").Value = 44593
$ws.Range("B115").Value = 1.5
$ws.Range("C115").Value = "Perusasetusten teko ja jest kirjastojen asennus, ekan testin rakentamista ja errorien korjausta"
$ws.Range("D115").Value = "client"

# --- new entry: row 116 (2.2.2022) --------------------------------------
$ws.Range("A113").Copy() | Out-Null
$ws.Range("A116").PasteSpecial(-4122) | Out-Null
$ws.Range("A116").Value = 44594
$ws.Range("B116").Value = 1
$ws.Range("C116").Value = "PrivateRoute testit"
$ws.Range("D116").Value = "client"

# --- totals block now lives at rows 119-121; extend the hours sum so it
# covers the three new rows too (percentage formula already shifted and
# recalculates automatically) -------------------------------------------
$ws.Range("B119").Formula = "=SUM(B2:B116)"

# --- view bookkeeping: match the saved selection/scroll position --------
$ws.Range("C117").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 1
